$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Distance-of-last-point-to-center-of-dest values in column C
$ws.Range("C1").Value = 0.75
$ws.Range("C2").Value = 0.5
$ws.Range("C3").Value = 0.25

# Touch cell protection/alignment formatting on the whole used range so a
# new (explicit) cell style gets minted for A1:C3, matching the format
# change captured upstream.
$ws.Range("A1:C3").Locked = $false

# Move the active selection to C1 (single cell) as in the edited workbook.
$ws.Range("C1").Select()
